$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("E2").Value = "2016-03-21 05:17:14"
$wsZh.Range("E3").Value = "2016-03-21 05:17:14"
$wsZh.Range("H2").Value = "2016-03-21 05:17:58"
$wsZh.Range("H3").Value = "2016-03-21 05:17:58"

$wsDe.Range("E2").Value = "2016-03-21 05:17:25"
$wsDe.Range("E3").Value = "2016-03-21 05:17:25"
$wsDe.Range("H2").Value = "2016-03-21 05:18:13"
$wsDe.Range("H3").Value = "2016-03-21 05:18:13"
